# Combine the two datasets (drop the now-redundant "and other" / split
# recycling-backfilling rows) and move the trailing spacer/link rows up to
# directly follow the remaining data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that duplicated/expanded the combined categories
# (RCV_R, RCV_B, and the old "Disposal - landfill and other" row), which
# shifts the trailing rows (old B21/B22) up to B18/B19 automatically.
$ws.Rows("7:9").Delete()

# Re-point the remaining treatment-code / full-name rows at the merged
# dataset's values.
$ws.Range("A2").Value = "DSP_L"
$ws.Range("B2").Value = "Disposal - landfill (D1, D5, D12)"
$ws.Range("C2").Value = "disposal"

$ws.Range("A3").Value = "DSP_I"
$ws.Range("B3").Value = "Disposal - incineration (D10)"
$ws.Range("C3").Value = "disposal"

$ws.Range("A4").Value = "DSP_OTH"
$ws.Range("B4").Value = "Disposal - other (D2-D4, D6-D7)"
$ws.Range("C4").Value = "disposal"

$ws.Range("A5").Value = "RCV_E"
$ws.Range("B5").Value = "Recovery - energy recovery (R1)"
$ws.Range("C5").Value = "recovery"

$ws.Range("A6").Value = "RCV_R_B"
$ws.Range("B6").Value = "Recovery - recycling and backfilling (R2-R11)"
$ws.Range("C6").Value = "recovery"

# Move the selection down to where the next batch of data will be entered.
$ws.Range("A7:XFD8").Select()
